$d = $word.ActiveDocument

function Find-ParaIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -match [regex]::Escape($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right after
#    the title heading: "Meta description: Enjoy Age of Halvar slot game
#    with Wild Halvar and bonus mode. Discover the pros and cons in our
#    review and play for free."
# ---------------------------------------------------------------------------
$metaIndex = Find-ParaIndex $d "Meta description"
if ($metaIndex -gt 0) {
    $metaPara = $d.Paragraphs.Item($metaIndex)
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Locate the last paragraph of the document (the leftover italic
#    "Create a feature image..." AI-prompt text) and, right before it,
#    insert a new bold paragraph reading:
#    "Play Age of Halvar Free: Review and Features"
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$targetIndex = $lastIndex - 1
$pBeforeLast = $d.Paragraphs.Item($targetIndex)

$insertionPoint = $pBeforeLast.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newIndex = $targetIndex + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.Style = "Normal"
$newPara.Range.Text = "Play Age of Halvar Free: Review and Features"

$newPara = $d.Paragraphs.Item($newIndex)
$newTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newTextRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) Replace the final paragraph's text (formerly the leftover image-prompt
#    text) with the meta-description copy, keeping its italic formatting.
# ---------------------------------------------------------------------------
$finalIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($finalIndex)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastTextRange.Text = "Enjoy Age of Halvar slot game with Wild Halvar and bonus mode. Discover the pros and cons in our review and play for free."
